$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 75, shifting rows 75-77 down to 76-78.
$ws.Rows.Item(75).Insert()

# Fill in the new row 75 with the weekly price entry.
$ws.Cells.Item(75, 1).Value = 10
$ws.Cells.Item(75, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(75, 3).Value = "La Araucanía"
$ws.Cells.Item(75, 4).Value = 44706
$ws.Cells.Item(75, 5).Value = 9
$ws.Cells.Item(75, 6).Value = 100112035
$ws.Cells.Item(75, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 50
$ws.Cells.Item(75, 11).Value = 30000
$ws.Cells.Item(75, 12).Value = 30000
$ws.Cells.Item(75, 13).Value = 30000
$ws.Cells.Item(75, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(75, 15).Value = "Región Metropolitana"
$ws.Cells.Item(75, 16).Value = 3000
$ws.Cells.Item(75, 17).Value = 10
$ws.Cells.Item(75, 18).Value = "Hortaliza"
